$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("a f f asdf a", 3, "teppy@teppy.com"),
    @("asdf asdf", 12, "te@ttp.com"),
    @("asdf ththth gh", 12, "asdf@adfas.cc"),
    @("asdf", 2, "asdf@fga.com"),
    @("saasd asdfasd", 12, "sadf@gds.com")
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
